# NIT-9003956632.xlsx "Estado de Cuenta" update
#
# The account-statement database was refreshed: the five "Periodo Mora"
# rows (B16:J20) that list period 1807..1803 (descending) are restated in
# ascending order (1803..1807). Everything else about the rows (worker,
# document id, amounts, formatting) stays the same - only the period
# label in column E is updated per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1803"
$ws.Range("E17").Value = "1804"
$ws.Range("E18").Value = "1805"
$ws.Range("E19").Value = "1806"
$ws.Range("E20").Value = "1807"
